# corrM.xlsx update: add two new variables ("skewMed" and "skewMean") to the
# correlation matrix on Sheet1.
#
# "skewMed" is inserted right after "iqrMed" (between the old columns F and G,
# and between the old rows 6 and 7) but no correlation data is available for
# it yet, so its entire row/column stays blank.
#
# "skewMean" is inserted right after "iqrMean" (between the old columns L/M
# and old rows 12/13) and DOES have correlation data, symmetric like the
# rest of the matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert the two new columns and two new rows.
#    Columns/rows shift right/down automatically, carrying values+styles.
# ---------------------------------------------------------------------

# New blank column for "skewMed" goes where old column G ("varMed") was.
$ws.Range("G1").EntireColumn.Insert()
# New column for "skewMean" goes where old column L+1 now sits, i.e.
# right after "iqrMean" (old column L) and before "varMean" (old column L,
# now shifted to M after the previous insert).
$ws.Range("M1").EntireColumn.Insert()

# New blank row for "skewMed" goes where old row 7 ("varMed") was.
$ws.Range("A7").EntireRow.Insert()
# New row for "skewMean" goes right after "iqrMean" (old row 12, now row 12
# still) and before "varMean" (old row 12 shifted to 13 after this insert).
$ws.Range("A13").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. Copy header/label formatting (bold, centered, bordered = same style
#    as every other header/label cell) onto the newly-inserted cells.
# ---------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the header row / label column text for the two new variables.
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "skewMed"
$ws.Range("A7").Value = "skewMed"

$ws.Range("M1").Value = "skewMean"
$ws.Range("A13").Value = "skewMean"

# ---------------------------------------------------------------------
# 4. Fill in the correlation values for "skewMean" (row 13 and column M).
#    Row/column G and row 7 (the "skewMed" data) are intentionally left
#    blank -- no data is available for that variable yet.
# ---------------------------------------------------------------------

# skewMean row (row 13), columns B..P except G (stays blank)
$ws.Range("B13").Value = -0.01711181371022921
$ws.Range("C13").Value = 0.01183533634683476
$ws.Range("D13").Value = -0.07888568017450877
$ws.Range("E13").Value = -0.04129872032365137
$ws.Range("F13").Value = -0.2322392139256818
$ws.Range("H13").Value = -0.1817930213045486
$ws.Range("I13").Value = -0.3032079420773576
$ws.Range("J13").Value = 0.117551294370681
$ws.Range("K13").Value = -0.06793267903473671
$ws.Range("L13").Value = 0.155101694151577
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 0.167881496477694
$ws.Range("O13").Value = -0.2465782624408502
$ws.Range("P13").Value = 0.0008191071168884024

# skewMean column (M), remaining rows (row 13 already set above, row 7 stays blank)
$ws.Range("M2").Value = -0.01711181371022921
$ws.Range("M3").Value = 0.01183533634683476
$ws.Range("M4").Value = -0.07888568017450877
$ws.Range("M5").Value = -0.04129872032365137
$ws.Range("M6").Value = -0.2322392139256818
$ws.Range("M8").Value = -0.1817930213045486
$ws.Range("M9").Value = -0.3032079420773576
$ws.Range("M10").Value = 0.117551294370681
$ws.Range("M11").Value = -0.06793267903473671
$ws.Range("M12").Value = 0.155101694151577
$ws.Range("M14").Value = 0.167881496477694
$ws.Range("M15").Value = -0.2465782624408502
$ws.Range("M16").Value = 0.0008191071168884024
